$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.64"
$ws.Range("D3").Value = "'23.54"
$ws.Range("D4").Value = "'5.640"
$ws.Range("D5").Value = "'0.05814"
$ws.Range("D6").Value = "'3.412"
$ws.Range("D7").Value = "'6.471"
$ws.Range("D9").Value = "'0.7973"
$ws.Range("D10").Value = "'0.1457"
$ws.Range("D12").Value = "'0.03237"
$ws.Range("D13").Value = "'0.02951"
$ws.Range("D14").Value = "'0.09233"
$ws.Range("D15").Value = "'0.001667"
$ws.Range("D16").Value = "'3.326"
$ws.Range("D17").Value = "'0.04750"
$ws.Range("D18").Value = "'0.0005996"
$ws.Range("D19").Value = "'0.006239"
$ws.Range("D20").Value = "'0.005461"
$ws.Range("D21").Value = "'0.001068"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.695"
$ws.Range("D26").Value = "'0.1236"
$ws.Range("D27").Value = "'0.001001"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("D40").Value = "'0.04304"
$ws.Range("D41").Value = "'0.007167"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003601"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1052"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.009526"
$ws.Range("D46").Value = "'0.00005757"
$ws.Range("D48").Value = "'0.7861"
$ws.Range("D49").Value = "'0.1065"
$ws.Range("E49").Value = "48BOLOBOLO"
